# WG N numbers table template
# Fill in the header row labels for the WG N# tracking table and move the
# active selection, matching the authored workbook.
#
# NB: values are written in the order they first appear in the shared
# string table of the authored file (Submitter, request date, Document and
# schema, N#) so the generated sharedStrings.xml lines up the same way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Submitter"
$ws.Range("F1").Value = "request date"
$ws.Range("E1").Value = "Document and schema"
$ws.Range("A1").Value = "N#"

# Restore the default zoom level (the source sheet had a stray 190% zoom).
$excel.ActiveWindow.Zoom = 100

# Move the cursor/selection to H13, like in the authored workbook.
[void]$ws.Range("H13").Select()
